$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "FT232500V3CNL7YV"
$ws.Range("A5").Value = "FT232500TL057FVV"
